$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '68.401.19'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  -3.61%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '3.700.91'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  -3.68%  '
$ws.Range('E4').Value = '  -0.25%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '584.42'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -1.48%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '183.45'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +10.49%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '3.693.02'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  -3.55%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.629'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -6.26%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.996'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -0.54%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.720'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -3.90%  '
$ws.Range('E11').Value = '  -6.57%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '56.33'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +6.15%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.0000292'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -8.65%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '10.44'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -6.26%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '4.191.69'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -6.19%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '3.688.75'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -4.52%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '19.44'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -5.83%  '
$ws.Range('E18').Value = '  -2.47%  '
$ws.Range('E19').Value = '  -6.12%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '12.81'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -7.13%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '68.150.38'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -3.93%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '410.00'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -5.24%  '
$ws.Range('E23').Value = '  -4.77%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '88.70'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -5.53%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '3.03'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -6.85%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '11.11'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +2.73%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '12.82'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -6.46%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '3.88'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -4.89%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '6.07'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +2.26%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '9.51'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -5.96%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '32.74'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -5.98%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '7.46'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -4.43%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '12.54'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -6.90%  '
$ws.Range('E34').Value = '  -5.71%  '
$ws.Range('B35').Value = 'InjectiveProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '43.88'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -12.90%  '
$ws.Range('B36').Value = 'OKB'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '65.28'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -4.96%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '594.29'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -3.31%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.0₃0893'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -8.66%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.403'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -3.25%  '
$ws.Range('E40').Value = '  +0.16%  '
$ws.Range('E41').Value = '  -0.36%  '
$ws.Range('E42').Value = '  -3.26%  '
$ws.Range('E43').Value = '  +5.67%  '
$ws.Range('E44').Value = '  -7.99%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '2.96'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -6.75%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.0436'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -6.23%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '9.32'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -7.97%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '2.778.84'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -1.60%  '
$ws.Range('E49').Value = '  -6.20%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '2.68'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -1.62%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '3.14'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -6.12%  '
